$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table is shrinking from 18 worker rows (16-33) to 14 (16-29). Row 33
# currently carries the "bottom of table" border formatting; copy that
# formatting onto row 29, which is the row that will become the new last
# row of the table once rows 30-33 are deleted below.
$ws.Range("B33:J33").Copy()
$ws.Range("B29:J29").PasteSpecial(-4122)

# --- Summary block updates ---
$ws.Range("E11").Value = 1151654
$ws.Range("C13").Value = 7
$ws.Range("F13").Value = 4

# --- Rebuilt worker/period table (rows 16-29) ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "80105195"
$ws.Range("D16").Value = "PEDRO IGNACIO ALVAREZ ACOSTA"
$ws.Range("E16").Value = "2505"
$ws.Range("F16").Value = 76000
$ws.Range("G16").Value = 1900000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1002320609"
$ws.Range("D17").Value = "LACIDES JOSE ESTRADA FLOREZ"
$ws.Range("E17").Value = "2505"
$ws.Range("F17").Value = 76000
$ws.Range("G17").Value = 1900000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1014736070"
$ws.Range("D18").Value = "JULIAN FELIPE AMAYA CARDENAS"
$ws.Range("E18").Value = "2505"
$ws.Range("F18").Value = 52000
$ws.Range("G18").Value = 1300000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "33332720"
$ws.Range("D19").Value = "SHIRLEY PAOLA MARTELO SANTOS"
$ws.Range("E19").Value = "2506"
$ws.Range("F19").Value = 88000
$ws.Range("G19").Value = 2200000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1010162749"
$ws.Range("D20").Value = "DIANA MARCELA CARDENAS RAMOS"
$ws.Range("E20").Value = "2506"
$ws.Range("F20").Value = 120000
$ws.Range("G20").Value = 3000000

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1014736070"
$ws.Range("D21").Value = "JULIAN FELIPE AMAYA CARDENAS"
$ws.Range("E21").Value = "2506"
$ws.Range("F21").Value = 52000
$ws.Range("G21").Value = 1300000

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "15646912"
$ws.Range("D22").Value = "JAIRO ENRIQUE GONZALEZ HERRERA"
$ws.Range("E22").Value = "2507"
$ws.Range("F22").Value = 43654
$ws.Range("G22").Value = 1423500

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "73203178"
$ws.Range("D23").Value = "ANUAR JOSE PARDO MORENO"
$ws.Range("E23").Value = "2507"
$ws.Range("F23").Value = 124000
$ws.Range("G23").Value = 3100000

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "33332720"
$ws.Range("D24").Value = "SHIRLEY PAOLA MARTELO SANTOS"
$ws.Range("E24").Value = "2507"
$ws.Range("F24").Value = 88000
$ws.Range("G24").Value = 2200000

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1010162749"
$ws.Range("D25").Value = "DIANA MARCELA CARDENAS RAMOS"
$ws.Range("E25").Value = "2507"
$ws.Range("F25").Value = 120000
$ws.Range("G25").Value = 3000000

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1014736070"
$ws.Range("D26").Value = "JULIAN FELIPE AMAYA CARDENAS"
$ws.Range("E26").Value = "2507"
$ws.Range("F26").Value = 52000
$ws.Range("G26").Value = 1300000

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "33332720"
$ws.Range("D27").Value = "SHIRLEY PAOLA MARTELO SANTOS"
$ws.Range("E27").Value = "2508"
$ws.Range("F27").Value = 88000
$ws.Range("G27").Value = 2200000

$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "1010162749"
$ws.Range("D28").Value = "DIANA MARCELA CARDENAS RAMOS"
$ws.Range("E28").Value = "2508"
$ws.Range("F28").Value = 120000
$ws.Range("G28").Value = 3000000

$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "1014736070"
$ws.Range("D29").Value = "JULIAN FELIPE AMAYA CARDENAS"
$ws.Range("E29").Value = "2508"
$ws.Range("F29").Value = 52000
$ws.Range("G29").Value = 1300000

# Remove the now-superseded worker rows (30-33); this shifts the footer
# (signature) block up from rows 38-39 to rows 34-35 and updates the used
# range / merged cells automatically.
$ws.Rows("30:33").Delete()
